$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Fortnite"
$ws.Range("B8").Value = "ryzen 3 3200G"

$ws.Range("C8").Formula = "=TEXT(367.4,""0.0"")"
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("E8").Select()
